# Applies two changes to Schedule-Links-list.docx:
# 1. Removes bold (w:b / w:bCs) from the paragraph-mark run properties of the
#    empty paragraph that sits right after "Panel - Q & A" and right before
#    "Genetics / DNA".
# 2. Prefixes the "10:30-12:10" paragraph (after "Bringing It All Together...")
#    with two new runs: "Thurs, June 6" (regular) and " - " (bold), matching
#    the pattern used for the other day/time headings in the document.

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $searchText, $exact = $true) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        # Paragraph text includes trailing paragraph mark char(s); trim them.
        $t = $t.TrimEnd([char]13, [char]7)
        if ($exact) {
            if ($t -eq $searchText) { return $p }
        } else {
            if ($t -like $searchText) { return $p }
        }
    }
    return $null
}

# --- Change 1: strip bold from the empty paragraph before "Genetics / DNA" ---

$genetics = Get-ParagraphByText $d "Genetics / DNA"
if ($null -eq $genetics) { throw "Could not find 'Genetics / DNA' paragraph" }

$emptyBeforeGenetics = $genetics.Previous()
if ($null -eq $emptyBeforeGenetics) { throw "Could not find empty paragraph before 'Genetics / DNA'" }
if ($emptyBeforeGenetics.Range.Text.TrimEnd([char]13, [char]7) -ne "") {
    throw "Paragraph before 'Genetics / DNA' was not empty as expected"
}

$xmlNoBold = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="360C33DC" w14:textId="77777777" w:rsidR="005108F5" w:rsidRPr="005108F5" w:rsidRDefault="005108F5" w:rsidP="007E1DAD"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

[void]$emptyBeforeGenetics.Range.InsertXML($xmlNoBold)

# --- Change 2: add "Thurs, June 6 - " prefix runs before "10:30-12:10" ---

$timeRange = Get-ParagraphByText $d "10:30-12:10"
if ($null -eq $timeRange) { throw "Could not find '10:30-12:10' paragraph" }

$xmlWithPrefix = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="02A7A4DC" w14:textId="77777777" w:rsidR="00317C31" w:rsidRPr="005108F5" w:rsidRDefault="00317C31" w:rsidP="00317C31"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Thurs, June 6</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r><w:r w:rsidRPr="005108F5"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-CA"/></w:rPr><w:t>10:30-12:10</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

[void]$timeRange.Range.InsertXML($xmlWithPrefix)

Write-Host "Edits applied successfully."
